# Add a new price-history snapshot column ("2026-01-31 23:11:58") right
# before the "nom" / "url_produit" columns, carrying forward the most
# recent known price (the previous last column, CP) for every product row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "nom" currently lives in column 95 (CQ); insert a new blank column there,
# which pushes "nom" -> CR (96) and "url_produit" -> CS (97).
$nomCol = 95
$ws.Columns.Item($nomCol).Insert()

$newCol = $nomCol          # 95 -> CQ, the freshly inserted column
$lastPriceCol = $nomCol - 1 # 94 -> CP, the previous most-recent snapshot

# Header for the new timestamp column. Column.Insert() already carries the
# surrounding header formatting (bold/border/centered) onto the new column,
# so we only need to set the value.
$ws.Cells.Item(1, $newCol).Value = "2026-01-31 23:11:58"

# Last data row.
$lastRow = $ws.UsedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $prevValue = $ws.Cells.Item($r, $lastPriceCol).Value()
    if ($prevValue -ne $null -and $prevValue -ne "") {
        $ws.Cells.Item($r, $newCol).Value = $prevValue
    }
}
